$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The weekly price sheet gets a new observation for this period. Insert a
# fresh row at position 143 (pushing the existing rows 143:229 down to
# 144:230, extending the used range to A1:R230) and populate it with the
# new week's data.
$ws.Rows.Item(143).Insert()

$ws.Cells.Item(143, 1).Value  = 3
$ws.Cells.Item(143, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(143, 3).Value  = "Coquimbo"
$ws.Cells.Item(143, 4).Value  = 44596
$ws.Cells.Item(143, 5).Value  = 5
$ws.Cells.Item(143, 6).Value  = 100112001
$ws.Cells.Item(143, 7).Value  = "Berenjena"
$ws.Cells.Item(143, 8).Value  = "Sin especificar"
$ws.Cells.Item(143, 9).Value  = "Primera"
$ws.Cells.Item(143, 10).Value = 73
$ws.Cells.Item(143, 11).Value = 9500
$ws.Cells.Item(143, 12).Value = 10000
$ws.Cells.Item(143, 13).Value = 9760
$ws.Cells.Item(143, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(143, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(143, 16).Value = 163
$ws.Cells.Item(143, 17).Value = 60
$ws.Cells.Item(143, 18).Value = "Hortaliza"

# Ensure the new row's date cell inherits the same date/time style (s="2")
# used throughout column D, same as the Rows.Insert carried-down formatting.
$ws.Cells.Item(143, 4).NumberFormat = $ws.Cells.Item(144, 4).NumberFormat
